$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150. This shifts the existing rows 150-232
# down to 151-233 (preserving all of their values/styles), which matches
# the bulk of the diff (every row from 151..233 equals the previous row's
# prior content).
$ws.Rows("150:150").Insert()

# Populate the newly inserted (blank) row 150 with its new data. Columns
# A, B, C, E, F, G, I, N, O, Q, R keep the same content the row had before
# the edit; only D, H, J, K, L, M, P receive new values per the diff.
$ws.Range("A150").Value = 10
$ws.Range("B150").Value = "Vega Modelo de Temuco"
$ws.Range("C150").Value = "La Araucanía"
$ws.Range("D150").Value = 44460
$ws.Range("E150").Value = 9
$ws.Range("F150").Value = 100112032
$ws.Range("G150").Value = "Zapallo italiano"
$ws.Range("H150").Value = "Bola 8"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 30
$ws.Range("K150").Value = 17000
$ws.Range("L150").Value = 17000
$ws.Range("M150").Value = 17000
$ws.Range("N150").Value = "$/caja 60 unidades"
$ws.Range("O150").Value = "Región de Arica y Parinacota"
$ws.Range("P150").Value = 283
$ws.Range("Q150").Value = 60
$ws.Range("R150").Value = "Hortaliza"
